$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell B1 = 0
$ws.Range("B1").Value = 0

# Row 2: A2 = 0, B2 = "disconnected_elements" (plain string, default style)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the shared style (bold font, thin box border, centered horizontal, top vertical)
# on B1 first, then copy just the formatting over to A2 so both cells end up pointing
# at the exact same style record instead of each mutation minting its own cellXf.
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108  # xlCenter
$b1.VerticalAlignment = -4160    # xlTop
$b1.Borders.Weight = 2           # xlThin

$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Save()
